$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update column header B2: "UFV_1" -> "UFV1"
$ws.Range("B2").Value = "UFV1"

# 2. Replace the "-" placeholder text in the data table with numeric 0
$dashCells = @(
    "B3","C3","D3",
    "B8","C8","D8",
    "B13","C13",
    "B15","C15","D15","E15",
    "B16","C16","D16",
    "B19","C19",
    "B25","C25","D25"
)
foreach ($addr in $dashCells) {
    $ws.Range($addr).Value = 0
}

# 3. Move the small legend table (I2:J7) down below the main data table (A29:B34)
$ws.Range("I2:J7").Cut($ws.Range("A29:B34"))
$ws.Range("I2:J7").Clear()
